$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Argentina Liga Profesional")

# Clear the old tail of the header row (AD1:AP1) completely (value + style),
# shrinking the sheet's used range/dimension back down.
$ws.Range("AD1:AP1").Clear()

# From Q1 onward the header labels shift left because the "_op" duplicate
# columns (Ah_op, oddAHH_op, AhOU_op, oddAHOver_op, oddAHUnder_op, PLH_op,
# PLD_op, PLA_op, PL_Ahh_op, PL_Aha_op, PL_AhOver_op, PL_AhUnder_op) were
# removed from the source data. Columns B1:P1 (id .. oddA) are untouched.
$ws.Range("Q1").Value  = "Ah"
$ws.Range("R1").Value  = "oddAHH"
$ws.Range("S1").Value  = "oddAHA"
$ws.Range("T1").Value  = "AhOU"
$ws.Range("U1").Value  = "oddAHOver"
$ws.Range("V1").Value  = "oddAHUnder"
$ws.Range("W1").Value  = "PLH"
$ws.Range("X1").Value  = "PLD"
$ws.Range("Y1").Value  = "PLA"
$ws.Range("Z1").Value  = "PL_Ahh"
$ws.Range("AA1").Value = "PL_Aha"
$ws.Range("AB1").Value = "PL_AhOver"
$ws.Range("AC1").Value = "PL_AhUnder"
